# "some changes to recovery bulkhead"
# Update the Black Powder Calculator worksheet with new measured inputs:
#   - Body tube length (B5): 13 -> 16 in
#   - Number of screws (B10): 4 -> 3
#   - Black powder weight (B11): 1.77 -> 3.2 g
# All downstream formulas (B15, B17-B20, B22) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Black Powder Calculator")

$ws.Range("B5").Value = 16
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 3.2

# Leave the cursor where the author left it when they saved the file.
$ws.Activate()
$null = $ws.Range("B19").Select()
